# Update NATMI LR-pair metrics (Tnfsf13-Tnfrsf14) with recomputed TPM-derived values.
# Source data changed for two clusters (ECs ligand avg/total, ECs receptor avg/total),
# which cascades into specificity (I/J/O/P) and edge weight (Q/R/S/T) columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5347833333333334
$ws.Range("H2").Value = 1.60435
$ws.Range("I2").Value = 0.196822066153855
$ws.Range("J2").Value = 0.196822066153855
$ws.Range("M2").Value = 3.961421333333333
$ws.Range("N2").Value = 11.884264
$ws.Range("O2").Value = 0.3114993985605504
$ws.Range("P2").Value = 0.3114993985605504
$ws.Range("Q2").Value = 2.118502105377778
$ws.Range("R2").Value = 19.0665189484
$ws.Range("S2").Value = 0.06130995523037071
$ws.Range("T2").Value = 0.06130995523037069
$ws.Range("G3").Value = 0.5347833333333334
$ws.Range("H3").Value = 1.60435
$ws.Range("I3").Value = 0.196822066153855
$ws.Range("J3").Value = 0.196822066153855
$ws.Range("O3").Value = 0.5009735319462221
$ws.Range("P3").Value = 0.500973531946222
$ws.Range("Q3").Value = 3.407112460155556
$ws.Range("R3").Value = 30.6640121414
$ws.Range("S3").Value = 0.09860264564604973
$ws.Range("T3").Value = 0.09860264564604969
$ws.Range("G4").Value = 0.5347833333333334
$ws.Range("H4").Value = 1.60435
$ws.Range("I4").Value = 0.196822066153855
$ws.Range("J4").Value = 0.196822066153855
$ws.Range("O4").Value = 0.1875270694932276
$ws.Range("P4").Value = 0.1875270694932276
$ws.Range("Q4").Value = 1.2753684064
$ws.Range("R4").Value = 11.4783156576
$ws.Range("S4").Value = 0.03690946527743462
$ws.Range("T4").Value = 0.03690946527743461
$ws.Range("I5").Value = 0.1891972429821067
$ws.Range("J5").Value = 0.1891972429821067
$ws.Range("M5").Value = 3.961421333333333
$ws.Range("N5").Value = 11.884264
$ws.Range("O5").Value = 0.3114993985605504
$ws.Range("P5").Value = 0.3114993985605504
$ws.Range("Q5").Value = 2.036432019141333
$ws.Range("R5").Value = 18.327888172272
$ws.Range("S5").Value = 0.05893482739824056
$ws.Range("T5").Value = 0.05893482739824054
$ws.Range("I6").Value = 0.1891972429821067
$ws.Range("J6").Value = 0.1891972429821067
$ws.Range("O6").Value = 0.5009735319462221
$ws.Range("P6").Value = 0.500973531946222
$ws.Range("S6").Value = 0.09478281105123358
$ws.Range("T6").Value = 0.09478281105123354
$ws.Range("I7").Value = 0.1891972429821067
$ws.Range("J7").Value = 0.1891972429821067
$ws.Range("O7").Value = 0.1875270694932276
$ws.Range("P7").Value = 0.1875270694932276
$ws.Range("S7").Value = 0.0354796045326326
$ws.Range("T7").Value = 0.0354796045326326
$ws.Range("I8").Value = 0.6139806908640383
$ws.Range("J8").Value = 0.6139806908640382
$ws.Range("M8").Value = 3.961421333333333
$ws.Range("N8").Value = 11.884264
$ws.Range("O8").Value = 0.3114993985605504
$ws.Range("P8").Value = 0.3114993985605504
$ws.Range("Q8").Value = 6.608605486541334
$ws.Range("R8").Value = 59.477449378872
$ws.Range("S8").Value = 0.1912546159319392
$ws.Range("T8").Value = 0.1912546159319391
$ws.Range("I9").Value = 0.6139806908640383
$ws.Range("J9").Value = 0.6139806908640382
$ws.Range("O9").Value = 0.5009735319462221
$ws.Range("P9").Value = 0.500973531946222
$ws.Range("S9").Value = 0.3075880752489388
$ws.Range("T9").Value = 0.3075880752489387
$ws.Range("I10").Value = 0.6139806908640383
$ws.Range("J10").Value = 0.6139806908640382
$ws.Range("O10").Value = 0.1875270694932276
$ws.Range("P10").Value = 0.1875270694932276
$ws.Range("S10").Value = 0.1151379996831604
$ws.Range("T10").Value = 0.1151379996831604
